# Update countries & provincias Spain
#
# 1. Refresh the "last updated" timestamp shown at the top of the sheet.
# 2. Refresh the covid stats for "Corea del Sur" (no reordering required).
# 3. Refresh the covid stats for "Mexico" and "Guatemala" with newer,
#    larger totals that now outrank their former neighbours
#    ("Arabia Saudita" and "Guadalupe" respectively).
# 4. Re-sort the country table (rows 4-216) by "Casos totales" (column B)
#    descending, which is how this sheet is always ordered, so the two
#    countries that grew naturally move above the countries they passed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Datos actualizados" timestamp (row 1) -----------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 03:52"

# --- 2. Update Corea del Sur's numbers -------------------------------------
$row = $ws.Cells.Find("Corea del Sur").Row
$ws.Cells.Item($row, 2).Value = 10512
$ws.Cells.Item($row, 3).Value = 32
$ws.Cells.Item($row, 4).Value = 7368
$ws.Cells.Item($row, 5).Value = 2930
$ws.Cells.Item($row, 6).Value = 55
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = 214

# --- 3. Update Mexico's numbers ---------------------------------------------
$row = $ws.Cells.Find("Mexico").Row
$ws.Cells.Item($row, 2).Value = 4219
$ws.Cells.Item($row, 3).Value = 375
$ws.Cells.Item($row, 4).Value = 1772
$ws.Cells.Item($row, 5).Value = 2174
$ws.Cells.Item($row, 6).Value = 89
$ws.Cells.Item($row, 7).Value = 40
$ws.Cells.Item($row, 8).Value = 273

# --- Update Guatemala's numbers ---------------------------------------------
$row = $ws.Cells.Find("Guatemala").Row
$ws.Cells.Item($row, 2).Value = 153
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 19
$ws.Cells.Item($row, 5).Value = 131
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 3

# --- 4. Re-sort the data table by Casos totales (column B), descending -----
$dataRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$dataRange.Sort($sortKey, 2)
